$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 10; I = "aa"; J = "Agree/Accept" },
    @{ Row = 14; I = "aa"; J = "Agree/Accept" },
    @{ Row = 16; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 21; I = "ba"; J = "Appreciation" },
    @{ Row = 23; I = "%"; J = "Uninterpretable" },
    @{ Row = 26; I = "qy"; J = "Yes-No-Question" },
    @{ Row = 37; I = "sv"; J = "Statement-opinion" },
    @{ Row = 50; I = "sv"; J = "Statement-opinion" },
    @{ Row = 55; I = "ba"; J = "Appreciation" },
    @{ Row = 62; I = "ba"; J = "Appreciation" },
    @{ Row = 67; I = "sv"; J = "Statement-opinion" },
    @{ Row = 71; I = "sv"; J = "Statement-opinion" },
    @{ Row = 76; I = "sv"; J = "Statement-opinion" },
    @{ Row = 77; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.I
    $ws.Range("J$($u.Row)").Value = $u.J
}
